$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("U4").Value = -3.736
$ws.Range("V4").Value = -1.817
$ws.Range("U7").Value = 70.7473
$ws.Range("V7").Value = -105.4
$ws.Range("U8").Value = -4.8321
$ws.Range("V8").Value = -5.7309
$ws.Range("U9").Value = -77.005
$ws.Range("V9").Value = -153.582
$ws.Range("U10").Value = 26.818
$ws.Range("V10").Value = -49.5337
$ws.Range("U11").Value = -2.1697
$ws.Range("V11").Value = -7.4158
$ws.Range("U12").Value = 68.7209
$ws.Range("V12").Value = -127.65
$ws.Range("U13").Value = -19.7602
$ws.Range("V13").Value = -21.3967
$ws.Range("U14").Value = 0.0663
$ws.Range("V14").Value = -6.3182
$ws.Range("U15").Value = -1.177
$ws.Range("V15").Value = -9.7224
$ws.Range("U17").Value = -8.3057
$ws.Range("V17").Value = -13.1546
$ws.Range("U18").Value = -17.3165
$ws.Range("V18").Value = -37.9019
$ws.Range("U19").Value = -11.2989
$ws.Range("V19").Value = -53.408
$ws.Range("U20").Value = 3.8551
$ws.Range("V20").Value = -23.077
$ws.Range("U21").Value = -1.967
$ws.Range("V21").Value = 0.3854
$ws.Range("U22").Value = -2.299
$ws.Range("V22").Value = -0.647
$ws.Range("U23").Value = 0.3518
$ws.Range("V23").Value = -2.795
$ws.Range("U24").Value = -1.088
$ws.Range("V24").Value = -1.385
$ws.Range("U25").Value = -2.206
$ws.Range("V25").Value = -1.154
$ws.Range("U26").Value = -3.255
$ws.Range("V26").Value = -1.148
$ws.Range("U27").Value = -2.804
$ws.Range("V27").Value = -1.614
$ws.Range("U28").Value = 0.8051
$ws.Range("V28").Value = -3.639
$ws.Range("U29").Value = -0.351
$ws.Range("V29").Value = -2.249
$ws.Range("U30").Value = 0.0589
$ws.Range("V30").Value = -1.712
$ws.Range("U31").Value = -4.437
$ws.Range("V31").Value = 0.6251
$ws.Range("U32").Value = -9.926
$ws.Range("V32").Value = -2.764
$ws.Range("U33").Value = -5.198
$ws.Range("V33").Value = 0.1474
$ws.Range("U34").Value = -6.12
$ws.Range("V34").Value = 0.972
$ws.Range("U35").Value = -0.908
$ws.Range("V35").Value = -2.547
$ws.Range("U37").Value = -4.084
$ws.Range("V37").Value = -0.328
$ws.Range("U39").Value = -4.027
$ws.Range("V39").Value = 0.0135
$ws.Range("U42").Value = -6.362
$ws.Range("V42").Value = 0.6562
$ws.Range("U45").Value = -7.808
$ws.Range("V45").Value = 0.2855
$ws.Range("U47").Value = -5.6444
$ws.Range("V47").Value = 0.6095
$ws.Range("U48").Value = -3.485
$ws.Range("V48").Value = -5.254
$ws.Range("U49").Value = -4.821
$ws.Range("V49").Value = -3.46
$ws.Range("U50").Value = 0.3575
$ws.Range("V50").Value = -6.892
$ws.Range("U52").Value = -2.492
$ws.Range("V52").Value = -4.503
$ws.Range("U54").Value = -5.2189
$ws.Range("V54").Value = -5.2462
$ws.Range("U56").Value = -5.027
$ws.Range("V56").Value = -3.025
$ws.Range("U57").Value = -3.923
$ws.Range("V57").Value = -2.987
$ws.Range("U62").Value = -1.7
$ws.Range("V62").Value = -0.986
$ws.Range("U63").Value = 0.5228
$ws.Range("V63").Value = -2.924
$ws.Range("U64").Value = -0.62
$ws.Range("V64").Value = -1.626
$ws.Range("U65").Value = -2.563
$ws.Range("V65").Value = -1.541
$ws.Range("U75").Value = -1.3593
$ws.Range("V75").Value = -3.7091
